$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.7287194209349384
$ws.Range("C2").Value = 1.65323645889881
$ws.Range("D2").Value = 2938.103010863317
$ws.Range("E2").Value = 71517.89157740913
$ws.Range("G2").Value = 74458.37654415228
